$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Row 3: update period-to-expire and last-update date
$ws.Range("H3").Value = -97
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "04-Nov-2025"

# Row 4: update period-to-expire and last-update date
$ws.Range("H4").Value = 699
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "04-Nov-2025"
